# Generate Report for Handback
# Updates the "Latest HO Xliff Generate Date" / "Correspond Handoff Datetime" /
# "Correspond Handback DateTime" timestamps for the c5ebdeb9-... row once its
# handback xliffs have actually been generated/processed.

$wb = $excel.ActiveWorkbook

# --- Overview sheet --------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G3").Value = "2016-09-01 04:52:48"

# --- zh-cn sheet -------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H3").Value = "2016-09-01 04:52:43"
$wsZhCn.Range("K3").Value = "2016-09-01 04:53:03"

# --- de-de sheet -------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H3").Value = "2016-09-01 04:52:48"
$wsDeDe.Range("K3").Value = "2016-09-01 04:53:14"
